$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-24 (Job ID, Company, Job Title, Candidate, Status, Action Date[serial])
$data = @(
    @(196, "Blockaid", "Enterprise Account Executive (Fintech)", "Brendan McMenimen", "2nd Interview", 45993),
    @(274, "Blockaid", "Senior Sales Engineer (US)", "Silja Petasch", "1st Interview", 45974),
    @(462, "Blockaid", "CS2 - Blockaid - Technical Account Manager (New York)", "Noah Newfield", "1st Interview", 45964),
    @(462, "Blockaid", "CS2 - Blockaid - Technical Account Manager (New York)", "Olivia Lo", "1st Interview", 45972),
    @(462, "Blockaid", "CS2 - Blockaid - Technical Account Manager (New York)", "Sean MacWilliams", "3rd Interview", 45971),
    @(541, "Blockaid", "CSM - Singapore", "Lim Yi Jun", "2nd Interview", 45995),
    @(669, "Blockaid", "Head of Marketing", "Jeff Meisel", "4th Interview", 45979),
    @(669, "Blockaid", "Head of Marketing", "Iulia Mihailescu", "1st Interview", 45985),
    @(766, "Cogent Security", "Enterprise Account Executive (US)", "Kevin Burns", "CV Sent", 45988),
    @(766, "Cogent Security", "Enterprise Account Executive (US)", "Jared Seavey", "1st Interview", 45980),
    @(766, "Cogent Security", "Enterprise Account Executive (US)", "Brendan McMenimen", "CV Sent", 45980),
    @(777, "Adaptive6", "SE Director", "Brendan Cox", "1st Interview", 45978),
    @(777, "Adaptive6", "SE Director", "Sean Valois", "2nd Interview", 45978),
    @(777, "Adaptive6", "SE Director", "Itai Heller", "4th Interview", 45989),
    @(777, "Adaptive6", "SE Director", "Hermann Hesse", "CV Sent", 45988),
    @(787, "Allium", "Solutions Engineer", "Silja Petasch", "3rd Interview", 45989),
    @(791, "Adaptive6", "Head of Sales (US)", "Brian Henger", "CV Sent", 45987),
    @(817, "Oasis Security", "Senior Sales Engineer (West/Mountain)", "Mary Greenlee", "1st Interview", 45985),
    @(824, "Blockaid", "Technical Account Manager", "Adam Palmer", "1st Interview", 45992),
    @(824, "Blockaid", "Technical Account Manager", "Derek Appleby", "CV Sent", 45986),
    @(826, "Legit Security", "VP of Sales", "Jim Underwood", "1st Interview", 45987),
    @(834, "Blockaid", "Regional Director US", "Greg Hilsenrath", "1st Interview", 45996),
    @(834, "Blockaid", "Regional Director US", "Jon Webster", "CV Sent", 45987)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
